$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# Re-ran the averaging code including the new spiral rotation schemes. This
# re-orders the "Gaussian-Quadrature" row right after the single/ring schemes,
# inserts three new "Spiral-*" scheme rows with freshly computed intensities,
# and pushes the remaining NoRotation/Rotation/HexGrid rows down accordingly.
# ----------------------------------------------------------------------------

# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9801883432228856
$ws.Range("D10").Value = 1.059928533930848
$ws.Range("E10").Value = 0.9808058266454883
$ws.Range("F10").Value = 0.9801883432228856
$ws.Range("G10").Value = 1.032208412466209
$ws.Range("H10").Value = 0.9500791628242369
$ws.Range("I10").Value = 0.9790772385206749
$ws.Range("J10").Value = 1.059928533930848
$ws.Range("K10").Value = 1.020367180288168
$ws.Range("L10").Value = 1.000277761755527
$ws.Range("M10").Value = 0.9970479196017238

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.003524112310321
$ws.Range("D11").Value = 0.9437060638926812
$ws.Range("E11").Value = 1.010503234934641
$ws.Range("F11").Value = 1.003524112310321
$ws.Range("G11").Value = 0.9655626310124884
$ws.Range("H11").Value = 1.034399559666349
$ws.Range("I11").Value = 1.009343399976948
$ws.Range("J11").Value = 0.9437060638926812
$ws.Range("K11").Value = 0.9771046494136613
$ws.Range("L11").Value = 0.9903143808619912
$ws.Range("M11").Value = 0.9945065002989047

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.003383937007914
$ws.Range("D12").Value = 0.9441006737053386
$ws.Range("E12").Value = 1.010408260722051
$ws.Range("F12").Value = 1.003383937007914
$ws.Range("G12").Value = 0.9657595508081134
$ws.Range("H12").Value = 1.034237759785192
$ws.Range("I12").Value = 1.009244342045347
$ws.Range("J12").Value = 0.9441006737053386
$ws.Range("K12").Value = 0.9772544672136947
$ws.Range("L12").Value = 0.990319202110804
$ws.Range("M12").Value = 0.9945224206789925

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.003481014248723
$ws.Range("D13").Value = 0.9438767101302591
$ws.Range("E13").Value = 1.010467632307549
$ws.Range("F13").Value = 1.003481014248723
$ws.Range("G13").Value = 0.9656233794518544
$ws.Range("H13").Value = 1.034381367540202
$ws.Range("I13").Value = 1.009307026785599
$ws.Range("J13").Value = 0.9438767101302591
$ws.Range("K13").Value = 0.9771721712189041
$ws.Range("L13").Value = 0.9903265927338135
$ws.Range("M13").Value = 0.9945228550773644

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9489560000000004
$ws.Range("D14").Value = 1.141036000000001
$ws.Range("E14").Value = 0.9595840000000003
$ws.Range("F14").Value = 0.9489560000000004
$ws.Range("G14").Value = 1.082087999999997
$ws.Range("H14").Value = 0.8887360000000017
$ws.Range("I14").Value = 0.9589560000000005
$ws.Range("J14").Value = 1.141036000000001
$ws.Range("K14").Value = 1.050310000000001
$ws.Range("L14").Value = 0.9996330000000005
$ws.Range("M14").Value = 0.9965593333333335

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.92
$ws.Range("D15").Value = 1.25
$ws.Range("E15").Value = 0.93
$ws.Range("F15").Value = 0.92
$ws.Range("G15").Value = 1.15
$ws.Range("H15").Value = 0.8
$ws.Range("I15").Value = 0.93
$ws.Range("J15").Value = 1.25
$ws.Range("K15").Value = 1.09
$ws.Range("L15").Value = 1.005
$ws.Range("M15").Value = 0.9966666666666666

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9511350392832012
$ws.Range("D16").Value = 1.144630985523201
$ws.Range("E16").Value = 0.9573382858752021
$ws.Range("F16").Value = 0.9511350392832012
$ws.Range("G16").Value = 1.084763097907189
$ws.Range("H16").Value = 0.8836034729984052
$ws.Range("I16").Value = 0.9574963456000006
$ws.Range("J16").Value = 1.144630985523201
$ws.Range("K16").Value = 1.050984635699201
$ws.Range("L16").Value = 1.001059837491201
$ws.Range("M16").Value = 0.9964945378645332

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Font.Bold = $true
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("A17").VerticalAlignment = -4160
$ws.Range("A17").Borders.LineStyle = 1
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.995803421894179
$ws.Range("D17").Value = 0.9953010795894998
$ws.Range("E17").Value = 0.9953074991545741
$ws.Range("F17").Value = 0.995803421894179
$ws.Range("G17").Value = 0.9958775954017556
$ws.Range("H17").Value = 0.9950023368591687
$ws.Range("I17").Value = 0.9954014828327942
$ws.Range("J17").Value = 0.9953010795894998
$ws.Range("K17").Value = 0.995304289372037
$ws.Range("L17").Value = 0.9955538556331079
$ws.Range("M17").Value = 0.9954489026219951

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Font.Bold = $true
$ws.Range("A18").HorizontalAlignment = -4108
$ws.Range("A18").VerticalAlignment = -4160
$ws.Range("A18").Borders.LineStyle = 1
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9985015000471089
$ws.Range("D18").Value = 0.985544179244658
$ws.Range("E18").Value = 0.9979527199521581
$ws.Range("F18").Value = 0.9985015000471089
$ws.Range("G18").Value = 0.9908930671588113
$ws.Range("H18").Value = 1.001199567704906
$ws.Range("I18").Value = 0.9977418191274302
$ws.Range("J18").Value = 0.985544179244658
$ws.Range("K18").Value = 0.991748449598408
$ws.Range("L18").Value = 0.9951249748227584
$ws.Range("M18").Value = 0.9953054755391788

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Style = $ws.Range("A16").Style
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.001818482498842
$ws.Range("D19").Value = 0.9663632798602781
$ws.Range("E19").Value = 1.002538476583146
$ws.Range("F19").Value = 1.001818482498842
$ws.Range("G19").Value = 0.9789405503306953
$ws.Range("H19").Value = 1.015319587876603
$ws.Range("I19").Value = 1.003237010886687
$ws.Range("J19").Value = 0.9663632798602781
$ws.Range("K19").Value = 0.984450878221712
$ws.Range("L19").Value = 0.993134680360277
$ws.Range("M19").Value = 0.9947028980060418

Write-Output ("Final UsedRange: " + $ws.UsedRange.Address())
